$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colIndexes = @(2, 3, 4, 6, 7, 9, 10, 12, 13)

$data = @(
    (@(22.57776839124094, 11.00185925728674, 4.273324926646968, 54.03919977346005, 3.775286636585458, 38.23205197220723, 9.798296835957828, 12.15427369529335, 19.79580128811288)),
    (@(22.36776654554976, 10.63651331433601, 4.227108077112563, 53.89259137818269, 3.779373839424425, 38.19468132075307, 9.80333678758136, 12.17413094065959, 19.78330701840359)),
    (@(22.24464504170485, 10.40968002296167, 4.198070356092988, 53.81353205062275, 3.782013154442323, 38.17842232907773, 9.806595350368942, 12.18794416562806, 19.7797575890323)),
    (@(22.19599119495421, 10.31678158731031, 4.186074150388775, 53.78408017967026, 3.783121451495909, 38.17347643286408, 9.807964571304243, 12.19398099471705, 19.77935026135742)),
    (@(22.18800557637348, 10.30133313396588, 4.184072425736872, 53.77935697500038, 3.783307465399331, 38.17275655630923, 9.808194428232115, 12.19500804461047, 19.7793454466314)),
    (@(22.24398265689241, 10.4084288009671, 4.197909227483644, 53.81312364466435, 3.782027968535021, 38.17834882878404, 9.806613648696921, 12.18802392881599, 19.77974788554045)),
    (@(22.50418598329445, 10.8765064960799, 4.25752818446958, 53.98638137022404, 3.776669049625039, 38.21777802450576, 9.800000623715954, 12.16078424353032, 19.79063942703988)),
    (@(23.05804308577931, 11.76769099945719, 4.369074637713533, 54.41254750993576, 3.767184008318156, 38.34819244061746, 9.788329793629773, 12.12021569415395, 19.84456914926637)),
    (@(23.4879121721897, 12.39795481554777, 4.447593645675741, 54.77743271113368, 3.760831394613848, 38.47637411370391, 9.780540322214458, 12.09822631450352, 19.90383197090238)),
    (@(23.68764666233932, 12.67788639612277, 4.482534870654556, 54.95442491466137, 3.758073467684892, 38.54169364542194, 9.777165940028153, 12.0899159923117, 19.93499848512526)),
    (@(23.76381362333104, 12.78280868688509, 4.495651586574382, 55.02300420414436, 3.757047947915024, 38.56743244465663, 9.775912375074869, 12.08701210200183, 19.94739947197113)),
    (@(23.74738716787041, 12.76026148091401, 4.492831817894624, 55.00816565930398, 3.757267975629515, 38.56184457148682, 9.776181275791286, 12.08762670333918, 19.94470216526378)),
    (@(23.69390270668225, 12.68654066579423, 4.48361630761045, 54.96003594488321, 3.757988720519763, 38.54379110388919, 9.777062323093217, 12.08967221869426, 19.93600675614287)),
    (@(23.66120911100125, 12.64124069488333, 4.477956508538493, 54.93075694541501, 3.758432648465873, 38.532863403773, 9.777605144631101, 12.09095679644253, 19.93075836953883)),
    (@(23.47493729402103, 12.37951480687421, 4.445294254321928, 54.76608510423588, 3.76101427561191, 38.47224603587981, 9.780764242127193, 12.0988034455586, 19.90187932192381)),
    (@(23.36168779045289, 12.21713833908755, 4.425055755513618, 54.66786485553345, 3.762631719022976, 38.43685179803431, 9.782745500294828, 12.1040504157742, 19.88523644437387)),
    (@(23.29694779787533, 12.12311126353658, 4.413342100600417, 54.61240865790541, 3.763574452796621, 38.41715401295772, 9.783900987269693, 12.10722768886667, 19.8760603430718)),
    (@(23.27509845183686, 12.09117014883812, 4.409363611239828, 54.59381108837027, 3.763895783626277, 38.4105981541684, 9.784294951805515, 12.10833083885848, 19.87302173909799)),
    (@(23.37370269439069, 12.23448985797768, 4.427217755818378, 54.67821337390579, 3.762458254545018, 38.44055128985802, 9.782532944895781, 12.10347537679755, 19.88696711944051)),
    (@(23.70959851403823, 12.70822440725107, 4.486326264409066, 54.97413079559338, 3.757776509726439, 38.54906664322365, 9.776802880853829, 12.08906480809939, 19.93854460638979)),
    (@(23.93219330974112, 13.01149085881005, 4.524286817208917, 55.17658665534909, 3.754826519476854, 38.62583538999103, 9.773199210427709, 12.0810631438345, 19.97574107898768)),
    (@(23.81313275915848, 12.85024483793997, 4.50408884641641, 55.06771276286775, 3.756390977907496, 38.58432906136949, 9.775109656292402, 12.08520430595195, 19.95557168818792)),
    (@(23.36826960274053, 12.22664732952414, 4.426240559055502, 54.67353165485199, 3.762536637820538, 38.438876722098, 9.782628989993935, 12.10373485128428, 19.88618345944738)),
    (@(22.90390949199011, 11.53035955533101, 4.339489828767221, 54.28808911341182, 3.769641197196187, 38.30722370106538, 9.791348779471889, 12.12981661840994, 19.82651255707456))
)

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $colIndexes.Count; $c++) {
        $ws.Cells.Item($excelRow, $colIndexes[$c]).Value2 = $rowValues[$c]
    }
}

$wb.Save()